{"js": "// Change the lottery draw time from \"19:25\" to \"18:25\" in the\n// preregistration text (\"... March 2025 at 19:25 pm as provided on the\n// Website ...\").\nconst body = context.document.body;\n\n// Find the run of text that contains the old time.\nconst results = body.search(\"19:25\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '19:25' in the document body.\");\n}\n\n// Replace just that found range's text with the corrected time.\nresults.items[0].insertText(\"18:25\", \"Replace\");\nawait context.sync();\n", "ps1": "# Change the lottery draw time from \"19:25\" to \"18:25\" in the\n# preregistration text (\"... March 2025 at 19:25 pm as provided on the\n# Website ...\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"19:25\"\n$find.Replacement.Text = \"18:25\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace:=wdReplaceAll\n$find.Execute(\"19:25\", $true, $false, $false, $false, $false, $true, 1, $false, \"18:25\", 2) | Out-Null\n"}
